$d = $word.ActiveDocument

# Disable smart-quote autocorrect so straight apostrophes in the
# replacement text are not converted to curly quotes.
$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false
$word.Options.AutoFormatReplaceQuotes = $false

function Replace-Text($old, $new) {
    $r = $d.Content
    $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $r.Text = $new
}

function Replace-All-Text($old, $new) {
    # Replace every occurrence of $old with $new (loops Find until no more
    # matches are found); used for strings like "[Muziki]" that occur more
    # than once in the document.
    $r = $d.Content
    while ($r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
        $r.Text = $new
        $r.Collapse(0)
        $r.End = $d.Content.End
    }
}

Replace-Text "Wafungwa na peremende - manukuu:" "Prisoners and candies - subtitles:"
Replace-Text "**mazungumzo huanza saa 55 ya pili sio 27 kwa sababu ya klipu ya utangulizi. Nilirekebisha nyakati ipasavyo. -John Argentino" "**dialogue starts at second 55 not 27 because of the intro clip. I adjusted the times accordingly. -John Argentino"

# Both [Muziki] occurrences become [Music]
Replace-All-Text "[Muziki]" "[Music]"

Replace-Text "wanahisabati wanne mkali wanachukuliwa" "four bright mathematicians are taken into"
Replace-Text "chini ya ulinzi na kuwekwa gerezani kwa sababu walijaribu" "custody and put in jail because they tried"
Replace-Text "kumshawishi mwanamke mzee kuwa Goedel's" "to convince an old lady that the Goedel's"
Replace-Text "nadharia za kutokamilika ni kweli. Kila" "incompleteness theorems are true. Every"
Replace-Text "mtaalamu wa hisabati ana kiini chake ambacho sisi" "mathematician has his own cell that we"
Replace-Text "inaweza kuhesabu na nambari kutoka 1 hadi 4." "can enumerate with a number from 1 to 4."
Replace-Text "kabla ya kuingia kwenye seli fulani" "before entering the cell a certain"
Replace-Text "idadi ya peremende kubwa kuliko " "number of candies greater than "
Replace-Text "e: AU " "e: OR "
Replace-Text "SAWA NA" "EQUAL TO"

# Delete the stray single-space run that sits between the ")" run and the
# " 1 ni"/" 1 is" run (the diff removes this run entirely).
$rParen = $d.Content
$rParen.Find.Execute(")", $false)
$rSpace = $d.Range($rParen.End, $rParen.End + 1)
$rSpace.Delete()

Replace-Text " 1 ni" " 1 is"
Replace-Text "wanapewa kila mtaalamu wa hisabati na wao" "given to every mathematician and they"
Replace-Text "wanaambiwa wana peremende 11 kwa jumla." "are told they have 11 candies in total."
Replace-Text "lakini kila mtu anajua idadi yake tu" "but everyone knows only his number of"
Replace-Text "pipi na jumla. 1 na sio" "candies and the total. 1 and is not"
Replace-Text "kuruhusiwa kuuliza nambari zingine." "allowed to ask for the others number."
Replace-Text "kisha mwanahisabati wa kwanza anauliza" "then the first mathematician asks the"
Replace-Text "pili: 'namba 2 unajua kama wewe" "second: 'number 2 do you know if you"
Replace-Text "kuwa na peremende nyingi kuliko mimi?' ya pili" "have more candies than me?' the second"
Replace-Text "mwanahisabati anajibu hana. Kisha yeye" "mathematician answers he doesn't. Then he"
Replace-Text "anauliza kwa nambari 3: 'unajua kama unayo" "asks to number 3: 'do you know if you have"
Replace-Text "pipi zaidi kuliko mimi?'" "more candy than me?'"
Replace-Text "mwanahisabati wa tatu anajibu: 'hapana niko" "the third mathematician answers: 'no I'm"
Replace-Text "samahani sifanyi'. Katika hatua hii ya nne" "sorry I don't'. At this point the fourth"
Replace-Text "mtaalamu wa hisabati anasema: 'jamani mnafahamu" "mathematician says: 'hey guys you know"
Replace-Text "nini, najua hasa pipi ngapi" "what, I know exactly how many candies"
Replace-Text "kila mtu ana hapa'. Cha kushangaza hata" "everyone has here'. Surprisingly even the"
Replace-Text "wanahisabati wengine watatu wanasema hivyo sasa" "other three mathematicians say that now"
Replace-Text "wanajua kila mtu ana pipi ngapi" "they know how many candies everyone has"
Replace-Text "kwa hivyo swali ni: unaweza kujua" "so the question is: can you figure out"
Replace-Text "idadi ya pipi kila mfungwa ana" "the number of candies every prisoner has"
